# Zero-pad the day-of-month "1" to "01" in several dates inside the first
# table, and append a new log row ("25/01/2021", "Allow 2 parameters...").
#
# Each date rewrite is expressed as a full-paragraph OOXML replacement so
# that the "0" that gets inserted lands in its own <w:r>, matching how Word
# splits runs when new text is typed in the middle of an existing run.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-ParagraphXml($range, $innerRunsXml) {
    $xml = $pkgHeader + '<w:body><w:p>' + $innerRunsXml + '</w:p></w:body>' + $pkgFooter
    $range.InsertXML($xml)
}

# Row index (1-based, including the header row) -> prefix kept before the
# split point, and whether that prefix is a separate pre-existing run.
# In every case the target cell text as a whole becomes:
#     <prefix> "0" "1/2021"
# where <prefix> ends in "/".

function Fix-DateCell($rowIndex, $prefixRunXml) {
    $cell = $t.Cell($rowIndex, 1)
    $runsXml = $prefixRunXml + '<w:r><w:t>0</w:t></w:r><w:r><w:t>1/2021</w:t></w:r>'
    Set-ParagraphXml $cell.Range $runsXml
}

# 1/1/2021 -> 1/ 0 1/2021
Fix-DateCell 7 '<w:r><w:t>1/</w:t></w:r>'

# 3/1/2021 -> 3/ 0 1/2021
Fix-DateCell 8 '<w:r><w:t>3/</w:t></w:r>'

# 7/1/2021 -> keep the pre-existing "7" run untouched, split "/1/2021"
Fix-DateCell 9 '<w:r><w:t>7</w:t></w:r><w:r><w:t>/</w:t></w:r>'

# 8/1/2021 -> keep the pre-existing "8" run (with its lastRenderedPageBreak)
Fix-DateCell 10 '<w:r><w:lastRenderedPageBreak/><w:t>8</w:t></w:r><w:r><w:t>/</w:t></w:r>'

# 10/1/2021 -> 10/ 0 1/2021
Fix-DateCell 11 '<w:r><w:t>10/</w:t></w:r>'

# 23/1/2021 -> 23/ 0 1/2021
Fix-DateCell 12 '<w:r><w:t>23/</w:t></w:r>'

# Append a brand-new row at the end of the table for the 25/01/2021 entry.
$newRow = $t.Rows.Add()

$dateXml = '<w:r><w:t>25/</w:t></w:r><w:r><w:t>0</w:t></w:r><w:r><w:t>1/2021</w:t></w:r>'
Set-ParagraphXml $newRow.Cells.Item(1).Range $dateXml

$descXml = '<w:r><w:t xml:space="preserve">Allow 2 parameters, branching factor for </w:t></w:r><w:r><w:t xml:space="preserve">existential and branching factor for universal </w:t></w:r>'
Set-ParagraphXml $newRow.Cells.Item(2).Range $descXml

# Columns 3 and 4 stay empty, matching the rest of the new row.
